$wb = $excel.ActiveWorkbook

# This script applies a scheduled-runner style data refresh to the market-price
# and profit columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values come from a re-scrape of current market prices; some rows gain/lose a
# cell in columns M/N depending on whether NQ/HQ profit is still negative.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1148.5
$ws.Range("I32").Value = 732.6667
$ws.Range("J32").Value = 1398
$ws.Range("K32").Value = 732.6667
$ws.Range("L32").Value = 1398
$ws.Range("M32").Value = -406.6667
$ws.Range("N32").Value = -2050
$ws.Range("H40").Value = 1076.6666
$ws.Range("I40").Value = 898.3333
$ws.Range("J40").Value = 1433.3334
$ws.Range("K40").Value = 898.3333
$ws.Range("L40").Value = 1433.3334
$ws.Range("M40").Value = -723.3333
$ws.Range("N40").Value = -1783.3334
$ws.Range("H76").Value = 4118637.2
$ws.Range("I76").Value = 4118637.2
$ws.Range("K76").Value = 4118637.2
$ws.Range("M76").Value = -4118322.2
$ws.Range("H79").Value = 4118637.2
$ws.Range("I79").Value = 4118637.2
$ws.Range("K79").Value = 4118637.2
$ws.Range("M79").Value = -4117545.2
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0
$ws.Range("H137").Value = 1776.6
$ws.Range("I137").Value = 1184
$ws.Range("J137").Value = 2261.4546
$ws.Range("K137").Value = 3552
$ws.Range("L137").Value = 6784.3638
$ws.Range("M137").Value = -1002
$ws.Range("N137").Value = -11884.3638
$ws.Range("H138").Value = 3405.1714
$ws.Range("I138").Value = 2449.138
$ws.Range("J138").Value = 4081.3901
$ws.Range("K138").Value = 7347.414
$ws.Range("L138").Value = 12244.1703
$ws.Range("M138").Value = -2207.414
$ws.Range("N138").Value = -22524.1703
$ws.Range("H139").Value = 78040
$ws.Range("J139").Value = 78040
$ws.Range("L139").Value = 78040
$ws.Range("N139").Value = -88320
$ws.Range("H140").Value = 96692.86
$ws.Range("J140").Value = 96692.86
$ws.Range("L140").Value = 96692.86
$ws.Range("N140").Value = -107052.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2306.2666
$ws.Range("I61").Value = 2507.2917
$ws.Range("J61").Value = 1502.1666
$ws.Range("K61").Value = 2507.2917
$ws.Range("L61").Value = 1502.1666
$ws.Range("M61").Value = -2295.2917
$ws.Range("N61").Value = -1926.1666
$ws.Range("H74").Value = 1126.2
$ws.Range("I74").Value = 1029.25
$ws.Range("J74").Value = 1514
$ws.Range("K74").Value = 1029.25
$ws.Range("L74").Value = 1514
$ws.Range("M74").Value = -155.25
$ws.Range("N74").Value = -3262
$ws.Range("H77").Value = 1126.2
$ws.Range("I77").Value = 1029.25
$ws.Range("J77").Value = 1514
$ws.Range("K77").Value = 5146.25
$ws.Range("L77").Value = 7570
$ws.Range("M77").Value = -778.25
$ws.Range("N77").Value = -16306
$ws.Range("H88").Value = 3885.4614
$ws.Range("I88").Value = 1801
$ws.Range("J88").Value = 5672.143
$ws.Range("K88").Value = 1801
$ws.Range("L88").Value = 5672.143
$ws.Range("M88").Value = -1395
$ws.Range("N88").Value = -6484.143
$ws.Range("H91").Value = 3885.4614
$ws.Range("I91").Value = 1801
$ws.Range("J91").Value = 5672.143
$ws.Range("K91").Value = 1801
$ws.Range("L91").Value = 5672.143
$ws.Range("M91").Value = -397
$ws.Range("N91").Value = -8480.143
$ws.Range("H132").Value = 2125.9058
$ws.Range("I132").Value = 1706.0465
$ws.Range("J132").Value = 3931.3
$ws.Range("K132").Value = 5118.139499999999
$ws.Range("L132").Value = 11793.9
$ws.Range("M132").Value = -2588.139499999999
$ws.Range("N132").Value = -16853.9
$ws.Range("H136").Value = 2306.2666
$ws.Range("I136").Value = 2507.2917
$ws.Range("J136").Value = 1502.1666
$ws.Range("K136").Value = 7521.875100000001
$ws.Range("L136").Value = 4506.4998
$ws.Range("M136").Value = -4971.875100000001
$ws.Range("N136").Value = -9606.4998
$ws.Range("H140").Value = 92532.55
$ws.Range("J140").Value = 92532.55
$ws.Range("L140").Value = 92532.55
$ws.Range("N140").Value = -102892.55
$ws.Range("H141").Value = 65571.73
$ws.Range("J141").Value = 65571.73
$ws.Range("L141").Value = 65571.73
$ws.Range("N141").Value = -75931.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 6995
$ws.Range("J19").Value = 6995
$ws.Range("L19").Value = 6995
$ws.Range("N19").Value = -7341
$ws.Range("H50").Value = 24326
$ws.Range("I50").Value = 100000
$ws.Range("J50").Value = 15917.777
$ws.Range("K50").Value = 100000
$ws.Range("L50").Value = 15917.777
$ws.Range("M50").Value = -99426
$ws.Range("N50").Value = -17065.777
$ws.Range("H134").Value = 1077.6428
$ws.Range("I134").Value = 1027.4
$ws.Range("J134").Value = 1203.25
$ws.Range("K134").Value = 3082.2
$ws.Range("L134").Value = 3609.75
$ws.Range("M134").Value = -547.2000000000003
$ws.Range("N134").Value = -8679.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12789.458
$ws.Range("I31").Value = 5789.5293
$ws.Range("J31").Value = 15622.762
$ws.Range("K31").Value = 5789.5293
$ws.Range("L31").Value = 15622.762
$ws.Range("M31").Value = -5494.5293
$ws.Range("N31").Value = -16212.762
$ws.Range("H34").Value = 12789.458
$ws.Range("I34").Value = 5789.5293
$ws.Range("J34").Value = 15622.762
$ws.Range("K34").Value = 5789.5293
$ws.Range("L34").Value = 15622.762
$ws.Range("M34").Value = -5587.5293
$ws.Range("N34").Value = -16026.762

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 100
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 100
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 1887.3
$ws.Range("I132").Value = 1245.5
$ws.Range("J132").Value = 2412.4092
$ws.Range("K132").Value = 11209.5
$ws.Range("L132").Value = 21711.6828
$ws.Range("M132").Value = -8679.5
$ws.Range("N132").Value = -26771.6828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5239806
$ws.Range("I113").Value = 10001346
$ws.Range("J113").Value = 911132.9399999999
$ws.Range("K113").Value = 10001346
$ws.Range("L113").Value = 911132.9399999999
$ws.Range("M113").Value = -9999176
$ws.Range("N113").Value = -915472.9399999999
$ws.Range("H138").Value = 69428.57000000001
$ws.Range("J138").Value = 69428.57000000001
$ws.Range("L138").Value = 69428.57000000001
$ws.Range("N138").Value = -79708.57000000001
$ws.Range("H141").Value = 45412
$ws.Range("J141").Value = 45412
$ws.Range("L141").Value = 45412
$ws.Range("N141").Value = -55772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10529358
$ws.Range("I40").Value = 3139.625
$ws.Range("K40").Value = 3139.625
$ws.Range("M40").Value = -3003.625
$ws.Range("H82").Value = 1826.7727
$ws.Range("I82").Value = 1977.1818
$ws.Range("J82").Value = 1676.3636
$ws.Range("K82").Value = 1977.1818
$ws.Range("L82").Value = 1676.3636
$ws.Range("M82").Value = -1616.1818
$ws.Range("N82").Value = -2398.3636
$ws.Range("H85").Value = 1826.7727
$ws.Range("I85").Value = 1977.1818
$ws.Range("J85").Value = 1676.3636
$ws.Range("K85").Value = 1977.1818
$ws.Range("L85").Value = 1676.3636
$ws.Range("M85").Value = -729.1818000000001
$ws.Range("N85").Value = -4172.3636
$ws.Range("H132").Value = 3450.9429
$ws.Range("I132").Value = 3675.6072
$ws.Range("J132").Value = 2552.2856
$ws.Range("K132").Value = 11026.8216
$ws.Range("L132").Value = 7656.8568
$ws.Range("M132").Value = -8496.821599999999
$ws.Range("N132").Value = -12716.8568
$ws.Range("H136").Value = 5662
$ws.Range("I136").Value = 4149.222
$ws.Range("J136").Value = 7023.5
$ws.Range("K136").Value = 12447.666
$ws.Range("L136").Value = 21070.5
$ws.Range("M136").Value = -9897.665999999999
$ws.Range("N136").Value = -26170.5
$ws.Range("H138").Value = 59168.816
$ws.Range("J138").Value = 59168.816
$ws.Range("L138").Value = 59168.816
$ws.Range("N138").Value = -69448.81599999999
$ws.Range("H140").Value = 54544.285
$ws.Range("I140").Value = 57960
$ws.Range("J140").Value = 53975
$ws.Range("K140").Value = 57960
$ws.Range("L140").Value = 53975
$ws.Range("M140").Value = -52780
$ws.Range("N140").Value = -64335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9933.333000000001
$ws.Range("I62").Value = 9933.333000000001
$ws.Range("K62").Value = 9933.333000000001
$ws.Range("M62").Value = -9309.333000000001
$ws.Range("H65").Value = 9933.333000000001
$ws.Range("I65").Value = 9933.333000000001
$ws.Range("K65").Value = 49666.665
$ws.Range("M65").Value = -46546.665
$ws.Range("H100").Value = 1140.6666
$ws.Range("I100").Value = 1245.375
$ws.Range("J100").Value = 303
$ws.Range("K100").Value = 2490.75
$ws.Range("L100").Value = 606
$ws.Range("M100").Value = -1949.75
$ws.Range("N100").Value = -1688
$ws.Range("H141").Value = 78916.25
$ws.Range("J141").Value = 78916.25
$ws.Range("L141").Value = 78916.25
$ws.Range("N141").Value = -89276.25
